# Apply: "solved climb stairs and house robber DP Problems"
#
# Adds two new tracker rows under the "Leetcode 150" section:
#   18: Climb Stairs  -> https://leetcode.com/problems/climbing-stairs/
#   19: House Robber  -> https://leetcode.com/problems/house-robber/description/
# and normalizes A17's date style (was the odd-one-out "style 3") to the
# same date style used everywhere else in column A ("style 1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right after row 17 (nothing below to shift, so this
# just appends 18:19) -- the new cells inherit row 17's own formatting
# (date style for col A, hyperlink style for col C) directly.
$ws.Rows("18:19").Insert(-4121)   # xlShiftDown

# --- Row 18: Climb Stairs (date serial 46080 = 2026-02-27) ---
$ws.Range("A18").Value2 = 46080
$ws.Range("B18").Value2 = "Climb Stairs"
$ws.Range("C18").Value2 = "https://leetcode.com/problems/climbing-stairs/"

# --- Row 19: House Robber (date serial 46080 = 2026-02-27) ---
$ws.Range("A19").Value2 = 46080
$ws.Range("B19").Value2 = "House Robber"
$ws.Range("C19").Value2 = "https://leetcode.com/problems/house-robber/description/"

# A17 and A18 both take the table's common date style ("style 1"); A19 keeps
# the style it inherited from the insert above.
$ws.Range("A15").Copy()
$ws.Range("A17").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A18").PasteSpecial(-4122)   # xlPasteFormats

# Register the two new hyperlinks (stamps the rIdNN relationships + the
# <hyperlink> refs for C18/C19).
$ws.Hyperlinks.Add($ws.Range("C18"), "https://leetcode.com/problems/climbing-stairs/")
$ws.Hyperlinks.Add($ws.Range("C19"), "https://leetcode.com/problems/house-robber/description/")

# Adding a hyperlink mints Excel's built-in "Hyperlink" named cell style, but
# this sheet already carries its own pre-existing hyperlink look (style "2",
# used by every other URL cell in the table) -- drop the newly-minted style...
$wb.Styles.Item("Hyperlink").Delete()

# ...and restore the sheet's existing hyperlink formatting (style "2", taken
# from C17) on the new URL cells (Delete() above resets their look).
$ws.Range("C17").Copy()
$ws.Range("C18:C19").PasteSpecial(-4122)   # xlPasteFormats
